# Update F column values (comment counts) across the three affected sheets
# to reflect newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 21437
$ws1.Range("F3").Value = 3383
$ws1.Range("F4").Value = 859
$ws1.Range("F7").Value = 812
$ws1.Range("F8").Value = 305
$ws1.Range("F12").Value = 594
$ws1.Range("F13").Value = 193
$ws1.Range("F14").Value = 374
$ws1.Range("F15").Value = 43
$ws1.Range("F16").Value = 472
$ws1.Range("F17").Value = 235
$ws1.Range("F18").Value = 47

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 97

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 21437
$ws4.Range("F7").Value = 3383
$ws4.Range("F8").Value = 859
$ws4.Range("F10").Value = 97
$ws4.Range("F13").Value = 812
$ws4.Range("F14").Value = 305
$ws4.Range("F23").Value = 594
$ws4.Range("F25").Value = 193
$ws4.Range("F27").Value = 374
$ws4.Range("F29").Value = 43
$ws4.Range("F30").Value = 472
$ws4.Range("F32").Value = 235
$ws4.Range("F33").Value = 47

$wb.Save()
